# Update (Removed Auto Arima)
# Updates the forecast figures on the "Forecast Comparison" sheet (columns
# C-G, rows 2-17) and the corresponding summary statistics on the "Summary"
# sheet (rows 9, 10, 11, 12, 14 in column B).

$wb = $excel.ActiveWorkbook

$forecastSheet = $wb.Worksheets.Item("Forecast Comparison")
$summarySheet  = $wb.Worksheets.Item("Summary")

# New values for columns C (Prophet Forecast), D (Amazon Mean Forecast),
# E (Amazon P70 Forecast), F (Amazon P80 Forecast), G (Amazon P90 Forecast)
# for each week row (2-17).
$newValues = @{
    2  = @(144, 340, 390, 432, 495)
    3  = @(193, 308, 360, 407, 478)
    4  = @(213, 291, 341, 385, 453)
    5  = @(193, 290, 340, 387, 458)
    6  = @(177, 296, 350, 403, 484)
    7  = @(211, 293, 345, 392, 463)
    8  = @(295, 286, 344, 404, 499)
    9  = @(388, 299, 357, 416, 508)
    10 = @(449, 291, 345, 397, 477)
    11 = @(478, 286, 341, 397, 485)
    12 = @(495, 289, 348, 411, 510)
    13 = @(503, 287, 349, 418, 529)
    14 = @(486, 291, 352, 419, 525)
    15 = @(457, 278, 337, 407, 517)
    16 = @(443, 267, 325, 395, 507)
    17 = @(462, 262, 319, 386, 493)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $forecastSheet.Cells.Item($row, 3).Value = $vals[0]
    $forecastSheet.Cells.Item($row, 4).Value = $vals[1]
    $forecastSheet.Cells.Item($row, 5).Value = $vals[2]
    $forecastSheet.Cells.Item($row, 6).Value = $vals[3]
    $forecastSheet.Cells.Item($row, 7).Value = $vals[4]
}

# Summary sheet - values are stored as text labels (not numbers) of the
# recomputed aggregates (Total Forecast 16/8/4 weeks, Max Forecast,
# Min Forecast). Force text storage (NumberFormat "@") so the digit
# strings aren't re-interpreted as numeric cells, then clear the
# number-format override so no stray style is left behind on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $summarySheet.Cells.Item(9, 2)  "5587"
Set-TextValue $summarySheet.Cells.Item(10, 2) "1814"
Set-TextValue $summarySheet.Cells.Item(11, 2) "743"
Set-TextValue $summarySheet.Cells.Item(12, 2) "503"
Set-TextValue $summarySheet.Cells.Item(14, 2) "144"
